$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Unprotect()

# Update the confidential disclaimer text (date 2021-03-18 -> 2021-03-19)
$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-19 for illustrative purposes only and are subject to change."

# Update the D (Weight) / E (Percent Change) values for rows 2-10
$ws.Range("D2").Value = 0.1063460595208126
$ws.Range("E2").Value = 0.024926686217009

$ws.Range("D3").Value = 0.1035217021513573
$ws.Range("E3").Value = 0.00935590631364569

$ws.Range("D4").Value = 0.1156809001984806
$ws.Range("E4").Value = -0.003496212436527202

$ws.Range("D5").Value = 0.1374883246703911
$ws.Range("E5").Value = 0.0006684491978610207

$ws.Range("D6").Value = 0.1319759454020601
$ws.Range("E6").Value = -0.002208155454144101

$ws.Range("D7").Value = 0.1466503385991087
$ws.Range("E7").Value = 0.00009676795045487907

$ws.Range("D8").Value = 0.1280903000199886
$ws.Range("E8").Value = 0.002275140300318501

$ws.Range("D9").Value = 0.1302464294378008
$ws.Range("E9").Value = 0.003355476314709449

$ws.Range("E10").Value = 0.00375808302300018
